$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.419.03"
$ws.Range("E2").Value = "  +0.82%  "

$ws.Range("D3").Value = "2.520.28"
$ws.Range("E3").Value = "  +2.25%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").Value = "  -0.20%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "519.72"
$ws.Range("E5").Value = "  +0.25%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "131.82"
$ws.Range("E6").Value = "  +0.82%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.08%  "

$ws.Range("E8").Value = "  -0.41%  "

$ws.Range("D9").Value = "2.518.35"
$ws.Range("E9").Value = "  +2.14%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0973"
$ws.Range("E10").Value = "  -1.78%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.21"
$ws.Range("E12").Value = "  -2.30%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.334"
$ws.Range("E13").Value = "  -2.30%  "

$ws.Range("D14").Value = "2.960.42"
$ws.Range("E14").Value = "  +2.05%  "

$ws.Range("D15").Value = "58.303.83"
$ws.Range("E15").Value = "  +0.79%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "22.21"
$ws.Range("E16").Value = "  -0.24%  "

$ws.Range("E17").Value = "  -0.69%  "

$ws.Range("D18").Value = "2.515.41"
$ws.Range("E18").Value = "  +2.26%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.72"
$ws.Range("E19").Value = "  -0.26%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "323.88"
$ws.Range("E20").Value = "  +1.22%  "

$ws.Range("E21").Value = "  +0.48%  "

$ws.Range("E22").Value = "  +5.74%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("E23").Value = "  -0.03%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.54"
$ws.Range("E24").Value = "  -0.75%  "

$ws.Range("E25").Value = "  -1.04%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.162"
$ws.Range("E26").Value = "  +1.02%  "

$ws.Range("E27").Value = "  -0.81%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.36"
$ws.Range("E28").Value = "  +0.72%  "

$ws.Range("D29").Value = "0.0₃0746"
$ws.Range("E29").Value = "  -0.71%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "168.02"
$ws.Range("E30").Value = "  +1.46%  "

$ws.Range("E31").Value = "  +0.56%  "

$ws.Range("B32").Value = "Aptos"
$ws.Range("C32").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.29"
$ws.Range("E32").Value = "  +0.04%  "

$ws.Range("B33").Value = "Fetch.AI"
$ws.Range("C33").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.18"
$ws.Range("E33").Value = "  +2.78%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.998"
$ws.Range("E34").Value = "  +0.00%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.998"
$ws.Range("E35").Value = "  -0.06%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.04"
$ws.Range("E36").Value = "  +0.06%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.27"
$ws.Range("E37").Value = "  -3.61%  "

$ws.Range("E38").Value = "  -0.77%  "

$ws.Range("E39").Value = "  +0.69%  "

$ws.Range("E40").Value = "  -1.19%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.778"
$ws.Range("E41").Value = "  -1.16%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "278.86"
$ws.Range("E42").Value = "  +2.60%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.10"
$ws.Range("E43").Value = "  +1.95%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.43"
$ws.Range("E44").Value = "  -0.43%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.598"
$ws.Range("E45").Value = "  +1.31%  "

$ws.Range("B46").Value = "Stellar"
$ws.Range("C46").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0921"
$ws.Range("E46").Value = "  +1.88%  "

$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "122.45"
$ws.Range("E47").Value = "  -2.60%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0498"
$ws.Range("E48").Value = "  +2.16%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "17.70"
$ws.Range("E49").Value = "  -0.51%  "

$ws.Range("E50").Value = "  +0.23%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "17.00"
$ws.Range("E51").Value = "  -0.20%  "
